$d = $word.ActiveDocument

# --- 1. Remove the "Meta description" paragraph that currently sits
#        right after the title (Heading1) paragraph:
#          "Meta description: Play Christmas Luck free and for real
#          money. Read our slot review and find recommended casinos to
#          try your holiday luck."
#        The whole paragraph (all its runs + its paragraph mark) is
#        deleted outright. ---
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# --- 2. The last paragraph in the document currently holds the italic
#        image-generation prompt ("Prompt: Create a cartoon-style
#        feature image..."). That paragraph is replaced by two
#        paragraphs: a new bold title line ("Play Christmas Luck Free:
#        Slot Review & Recommended Casinos") followed by the (still
#        italic) meta-description text that used to live at the top of
#        the document. ---
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)

# Target the paragraph's content only (exclude its trailing paragraph
# mark) so the replacement XML's own paragraph marks define the new
# paragraph boundaries, instead of an extra blank paragraph being
# minted.
$target = $d.Range($pLast.Range.Start, $pLast.Range.End - 1)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Christmas Luck Free: Slot Review &amp; Recommended Casinos</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Play Christmas Luck free and for real money. Read our slot review and find recommended casinos to try your holiday luck.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml)
